{"js": "// Add a new paragraph style \"FootnoteBlockText\" (\"Footnote Block Text\"),\n// based on the \"Footnote Text\" style, for block quotes inside footnotes.\n//\n// Passing the spaced display name to addStyle() makes the host derive the\n// styleId by stripping spaces (\"FootnoteBlockText\") while keeping the\n// w:name value intact (\"Footnote Block Text\") - matching how the style is\n// declared in styles.xml.\ncontext.document.addStyle(\"Footnote Block Text\", Word.StyleType.paragraph);\nawait context.sync();\n\n// Re-fetch the style by its id through the styles collection: the object\n// returned directly by addStyle() does not reliably resolve to the new\n// style for follow-up property writes, but a collection lookup does.\nconst style = context.document.getStyles().getByName(\"FootnoteBlockText\");\n\n// basedOn / next (the source document references the *display* name of\n// \"Footnote Text\", not its style id, so reproduce that literally).\nstyle.baseStyle = \"Footnote Text\";\nstyle.nextParagraphStyle = \"Footnote Text\";\n\nstyle.priority = 9;\nstyle.unhideWhenUsed = true;\nstyle.quickStyle = true;\nawait context.sync();\n\n// Paragraph formatting: spacing before/after = 100 twips (5pt) and\n// indentation left/right = 480 twips (24pt), first line = 0.\n// Word.ParagraphFormat properties are expressed in points, so divide the\n// twentieths-of-a-point (dxa) values from the target markup by 20.\nstyle.paragraphFormat.spaceBefore = 100 / 20;\nstyle.paragraphFormat.spaceAfter = 100 / 20;\nstyle.paragraphFormat.firstLineIndent = 0 / 20;\nstyle.paragraphFormat.leftIndent = 480 / 20;\nstyle.paragraphFormat.rightIndent = 480 / 20;\nawait context.sync();\n", "ps1": "# Add a new paragraph style \"FootnoteBlockText\" (\"Footnote Block Text\"),\n# based on the \"Footnote Text\" style, for block quotes inside footnotes.\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeParagraph = 1\n$style = $d.Styles.Add(\"FootnoteBlockText\", 1)\n\n# Word derives the style's NameLocal from the Add() name by default, so set\n# the spaced display name explicitly (the styleId stays \"FootnoteBlockText\").\n$style.NameLocal = \"Footnote Block Text\"\n\n# basedOn / next (the source document references the *display* name of\n# \"Footnote Text\", not its style id, so reproduce that literally).\n$style.BaseStyle = \"Footnote Text\"\n$style.NextParagraphStyle = \"Footnote Text\"\n\n$style.Priority = 9\n$style.UnhideWhenUsed = $true\n$style.QuickStyle = $true\n\n# Paragraph formatting: spacing before/after = 100 twips (5pt) and\n# indentation left/right = 480 twips (24pt), first line = 0.\n# Word's ParagraphFormat properties are expressed in points, so divide the\n# twentieths-of-a-point (dxa) values from the target markup by 20.\n$style.ParagraphFormat.SpaceBefore = 100 / 20\n$style.ParagraphFormat.SpaceAfter = 100 / 20\n$style.ParagraphFormat.FirstLineIndent = 0 / 20\n$style.ParagraphFormat.LeftIndent = 480 / 20\n$style.ParagraphFormat.RightIndent = 480 / 20\n"}
